$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) text updates ---
$ws.Range("A1").Value = "날짜"
$ws.Range("M1").Value = "엽면적지수"
$ws.Range("N1").Value = "주간생육길이_생육상태"
$ws.Range("O1").Value = "줄기굵기_생육상태"
$ws.Range("P1").Value = "잎길이_생육상태"
$ws.Range("Q1").Value = "입폭_생육상태"
$ws.Range("R1").Value = "입수_생육상태"
$ws.Range("S1").Value = "엽면적지수_생육상태"
$ws.Range("T1").Value = "개화화방위치_생육상태"
$ws.Range("U1").Value = "꽃과줄기거리_생육상태"
$ws.Range("V1").Value = "생육상태점수"
$ws.Range("W1").Value = "생장구분"

# --- Column A: convert "N주(MM/DD)" week labels into numeric dates 2018MMDD ---
# Each of the 15 weekly blocks spans exactly 16 data rows (B column sample numbers 1..16).
$weekDates = @(20180312, 20180319, 20180326, 20180402, 20180409, 20180416, 20180423, 20180430, 20180507, 20180514, 20180521, 20180528, 20180604, 20180611, 20180618)

$row = 2
foreach ($d in $weekDates) {
    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($row, 1).Value = $d
        $row++
    }
}

# --- growth_type_score (V) / growth_type (W) recalculated values for specific rows ---
$vwChanges = @{
    4 = @{ V = -4 }
    9 = @{ V = -6 }
    19 = @{ V = -6 }
    31 = @{ V = -1 }
    33 = @{ V = -1 }
    34 = @{ V = -5 }
    35 = @{ V = -5 }
    51 = @{ V = -4 }
    52 = @{ V = -4 }
    59 = @{ V = -5 }
    61 = @{ V = -5 }
    65 = @{ V = -5 }
    66 = @{ V = -4 }
    67 = @{ V = -3 }
    68 = @{ V = -2 }
    70 = @{ V = -4 }
    71 = @{ V = -5 }
    72 = @{ V = -4 }
    73 = @{ V = -5 }
    76 = @{ V = -5 }
    79 = @{ V = -5 }
    80 = @{ V = -4 }
    82 = @{ V = 0; W = 0 }
    83 = @{ V = 1; W = 1 }
    84 = @{ V = 0; W = 0 }
    85 = @{ V = -1 }
    86 = @{ V = -2 }
    87 = @{ V = -1 }
    88 = @{ V = -2 }
    89 = @{ V = -2 }
    90 = @{ V = -4 }
    91 = @{ V = -2 }
    92 = @{ V = -2 }
    93 = @{ V = 0; W = 0 }
    95 = @{ V = -2 }
    96 = @{ V = -1 }
    97 = @{ V = -2 }
    98 = @{ V = -2 }
    99 = @{ V = -1 }
    100 = @{ V = 0; W = 0 }
    101 = @{ V = -2 }
    102 = @{ V = -2 }
    103 = @{ V = -2 }
    104 = @{ V = -2 }
    105 = @{ V = 0; W = 0 }
    106 = @{ V = -2 }
    107 = @{ V = -2 }
    108 = @{ V = 0; W = 0 }
    109 = @{ V = -1 }
    110 = @{ V = -3 }
    111 = @{ V = -2 }
    112 = @{ V = -2 }
    113 = @{ V = 1; W = 1 }
    114 = @{ V = -5 }
    116 = @{ V = -4 }
    117 = @{ V = -5 }
    118 = @{ V = -3 }
    119 = @{ V = -5 }
    121 = @{ V = -4 }
    122 = @{ V = -5 }
    123 = @{ V = -5 }
    125 = @{ V = -4 }
    126 = @{ V = -3 }
    127 = @{ V = -2 }
    129 = @{ V = -2 }
    130 = @{ V = -4 }
    131 = @{ V = -3 }
    132 = @{ V = -2 }
    133 = @{ V = -4 }
    134 = @{ V = 0; W = 0 }
    135 = @{ V = -4 }
    137 = @{ V = -3 }
    138 = @{ V = -5 }
    139 = @{ V = -3 }
    140 = @{ V = -2 }
    141 = @{ V = -4 }
    143 = @{ V = 0; W = 0 }
    144 = @{ V = -1 }
    145 = @{ V = 0; W = 0 }
    147 = @{ V = -3 }
    148 = @{ V = -1 }
    150 = @{ V = 4 }
    151 = @{ V = -1 }
    152 = @{ V = -4 }
    154 = @{ V = -5 }
    156 = @{ V = 1; W = 1 }
    157 = @{ V = 0; W = 0 }
    158 = @{ V = -6 }
    159 = @{ V = 0; W = 0 }
    160 = @{ V = -3 }
    161 = @{ V = -2 }
    162 = @{ V = -3 }
    170 = @{ V = -2 }
    171 = @{ V = -1 }
    172 = @{ V = 4 }
    173 = @{ V = 2; W = 1 }
    174 = @{ V = -1 }
    175 = @{ V = -2 }
    176 = @{ V = 1; W = 1 }
    177 = @{ V = -3 }
    188 = @{ V = 2; W = 1 }
    189 = @{ V = 2; W = 1 }
    191 = @{ V = 1; W = 1 }
    192 = @{ V = 1; W = 1 }
    193 = @{ V = 2; W = 1 }
}

foreach ($r in $vwChanges.Keys) {
    $entry = $vwChanges[$r]
    $ws.Cells.Item($r, 22).Value = $entry.V
    if ($entry.ContainsKey("W")) {
        $ws.Cells.Item($r, 23).Value = $entry.W
    }
}

Write-Host "edit complete"
